$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sorel_format")

# Row 2 becomes the old "Lower Gorge (Columbia)" line (previously row 9).
$ws.Range("B2").Value = "Lower Gorge (Columbia)"
$ws.Range("C2").Value = 94
$ws.Range("D2").Value = "Gray_PC"

# Row 3 becomes a brand new "Kalama" line.
$ws.Range("B3").Value = "Kalama"
$ws.Range("C3").Value = 609
$ws.Range("D3").Value = "Gray_PC"

# Give the new numeric cells the same "0" integer format used by the old
# Lower Gorge row so they match style 10 (centered, numFmtId 1).
$ws.Range("C2:C3").NumberFormat = "0"
$ws.Range("C2:C3").HorizontalAlignment = -4108

# Remove the now-obsolete rows 4-9 (Coweeman, East Fork Lewis, Elochoman-
# Skamokawa, Grays-Chinook, South Fork Toutle, Green River, Washougal data
# that used to live in rows 2-8) and shift everything below up.
$ws.Range("A4:D9").Delete(-4162)

# Restore the view to match the saved workbook.
[void]$ws.Range("C10").Select()
